# Update TPM-derived values in the NATMI LR-pairs output sheet.
# The underlying receptor-expressing cell count for the "ECs" target
# cluster changed from 2 to 3 (with new TPM data), which cascades into
# several derived specificity / edge-weight columns for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs)
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.451416666666667
$ws.Range("N2").Value = 7.35425
$ws.Range("O2").Value = 0.2191928499183569
$ws.Range("P2").Value = 0.2191928499183569
$ws.Range("Q2").Value = 0.7876736776944445
$ws.Range("R2").Value = 7.089063099250001
$ws.Range("S2").Value = 0.1633085858138655
$ws.Range("T2").Value = 0.1633085858138655

# Row 3 (Target cluster: FAPs)
$ws.Range("O3").Value = 0.4446889938320204
$ws.Range("P3").Value = 0.4446889938320204
$ws.Range("S3").Value = 0.3313134107100094
$ws.Range("T3").Value = 0.3313134107100094

# Row 4 (Target cluster: MuSCs)
$ws.Range("O4").Value = 0.3361181562496228
$ws.Range("P4").Value = 0.3361181562496228
$ws.Range("S4").Value = 0.2504232267792271
$ws.Range("T4").Value = 0.2504232267792271

# Row 5 (Target cluster: ECs)
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.451416666666667
$ws.Range("N5").Value = 7.35425
$ws.Range("O5").Value = 0.2191928499183569
$ws.Range("P5").Value = 0.2191928499183569
$ws.Range("Q5").Value = 0.2695422510277778
$ws.Range("R5").Value = 2.42588025925
$ws.Range("S5").Value = 0.05588426410449135
$ws.Range("T5").Value = 0.05588426410449135

# Row 6 (Target cluster: FAPs)
$ws.Range("O6").Value = 0.4446889938320204
$ws.Range("P6").Value = 0.4446889938320204
$ws.Range("S6").Value = 0.113375583122011
$ws.Range("T6").Value = 0.113375583122011

# Row 7 (Target cluster: MuSCs)
$ws.Range("O7").Value = 0.3361181562496228
$ws.Range("P7").Value = 0.3361181562496228
$ws.Range("S7").Value = 0.08569492947039564
$ws.Range("T7").Value = 0.08569492947039564
